# Auto-generated: update Famfrit Profits market-data values per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 17.166666
$ws.Range("I6").Value = 17.166666
$ws.Range("K6").Value = 51.499998
$ws.Range("M6").Value = 60.500002

# Row 12
$ws.Range("H12").Value = 399.125
$ws.Range("I12").Value = 406.14285
$ws.Range("K12").Value = 406.14285
$ws.Range("M12").Value = -236.14285

# Row 29
$ws.Range("H29").Value = 2719.5
$ws.Range("I29").Value = 1079.25
$ws.Range("K29").Value = 3237.75
$ws.Range("M29").Value = -2956.75

# Row 58
$ws.Range("H58").Value = 519.6
$ws.Range("I58").Value = 519.6
$ws.Range("K58").Value = 1558.8
$ws.Range("M58").Value = -1408.8

# Row 138
$ws.Range("H138").Value = 22231816
$ws.Range("J138").Value = 30315386
$ws.Range("L138").Value = 90946158
$ws.Range("N138").Value = -90956438

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6768.5713
$ws.Range("I32").Value = 6426.6665
$ws.Range("J32").Value = 16000
$ws.Range("K32").Value = 6426.6665
$ws.Range("L32").Value = 16000
$ws.Range("M32").Value = -6139.6665
$ws.Range("N32").Value = -16574

# Row 61
$ws.Range("H61").Value = 11367973
$ws.Range("I61").Value = 12199166
$ws.Range("J61").Value = 8331.666999999999
$ws.Range("K61").Value = 12199166
$ws.Range("L61").Value = 8331.666999999999
$ws.Range("M61").Value = -12198954
$ws.Range("N61").Value = -8755.666999999999

# Row 136
$ws.Range("H136").Value = 11367973
$ws.Range("I136").Value = 12199166
$ws.Range("J136").Value = 8331.666999999999
$ws.Range("K136").Value = 36597498
$ws.Range("L136").Value = 24995.001
$ws.Range("M136").Value = -36594948
$ws.Range("N136").Value = -30095.001

# Row 138
$ws.Range("H138").Value = 23390
$ws.Range("I138").Value = 23390
$ws.Range("K138").Value = 23390
$ws.Range("M138").Value = -18250

# Row 139
$ws.Range("H139").Value = 134235.72
$ws.Range("J139").Value = 150000
$ws.Range("L139").Value = 150000
$ws.Range("N139").Value = -160280

# Row 141
$ws.Range("H141").Value = 69354.8
$ws.Range("J141").Value = 80846
$ws.Range("L141").Value = 80846
$ws.Range("N141").Value = -91206

$ws = $wb.Worksheets.Item("BSM")
# Row 138
$ws.Range("H138").Value = 199999
$ws.Range("J138").Value = 199999
$ws.Range("L138").Value = 199999
$ws.Range("N138").Value = -210279

# Row 140
$ws.Range("H140").Value = 119092.27
$ws.Range("J140").Value = 119092.27
$ws.Range("L140").Value = 119092.27
$ws.Range("N140").Value = -129452.27

# Row 141
$ws.Range("H141").Value = 66350
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3509.8472
$ws.Range("I31").Value = 2082.6538
$ws.Range("J31").Value = 7220.55
$ws.Range("K31").Value = 2082.6538
$ws.Range("L31").Value = 7220.55
$ws.Range("M31").Value = -1787.6538
$ws.Range("N31").Value = -7810.55

# Row 34
$ws.Range("H34").Value = 3509.8472
$ws.Range("I34").Value = 2082.6538
$ws.Range("J34").Value = 7220.55
$ws.Range("K34").Value = 2082.6538
$ws.Range("L34").Value = 7220.55
$ws.Range("M34").Value = -1880.6538
$ws.Range("N34").Value = -7624.55

# Row 69
$ws.Range("H69").Value = 25966
$ws.Range("I69").Value = 25966
$ws.Range("K69").Value = 25966
$ws.Range("M69").Value = -25217

# Row 72
$ws.Range("H72").Value = 25966
$ws.Range("I72").Value = 25966
$ws.Range("K72").Value = 77898
$ws.Range("M72").Value = -74154

# Row 132
$ws.Range("H132").Value = 46724.11
$ws.Range("I132").Value = 55497.21
$ws.Range("K132").Value = 166491.63
$ws.Range("M132").Value = -163961.63

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1632.7317
$ws.Range("J131").Value = 1757.3334
$ws.Range("L131").Value = 5272.0002
$ws.Range("N131").Value = -15352.0002

# Row 137
$ws.Range("H137").Value = 4206.5386
$ws.Range("I137").Value = 3849.5
$ws.Range("J137").Value = 4271.4546
$ws.Range("K137").Value = 11548.5
$ws.Range("L137").Value = 12814.3638
$ws.Range("M137").Value = -6448.5
$ws.Range("N137").Value = -23014.3638

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 1392555.2
$ws.Range("I14").Value = 1547269.8
$ws.Range("J14").Value = 125
$ws.Range("K14").Value = 1547269.8
$ws.Range("L14").Value = 125
$ws.Range("M14").Value = -1547101.8
$ws.Range("N14").Value = -461

# Row 102
$ws.Range("H102").Value = 3371.0833
$ws.Range("I102").Value = 1960.8334
$ws.Range("J102").Value = 4781.3335
$ws.Range("K102").Value = 1960.8334
$ws.Range("L102").Value = 4781.3335
$ws.Range("M102").Value = -338.8334
$ws.Range("N102").Value = -8025.3335

# Row 122
$ws.Range("H122").Value = 1435.9412
$ws.Range("I122").Value = 1424.5385
$ws.Range("J122").Value = 1473
$ws.Range("K122").Value = 4273.6155
$ws.Range("L122").Value = 4419
$ws.Range("M122").Value = -1823.6155
$ws.Range("N122").Value = -9319

# Row 123
$ws.Range("H123").Value = 25998.285
$ws.Range("J123").Value = 9997
$ws.Range("L123").Value = 9997
$ws.Range("N123").Value = -14897

# Row 126
$ws.Range("H126").Value = 7771.636
$ws.Range("I126").Value = 7926.857
$ws.Range("K126").Value = 23780.571
$ws.Range("M126").Value = -21310.571

# Row 132
$ws.Range("H132").Value = 3598.6667
$ws.Range("J132").Value = 2799.4
$ws.Range("L132").Value = 8398.200000000001
$ws.Range("N132").Value = -13458.2

# Row 135
$ws.Range("H135").Value = 199995
$ws.Range("J135").Value = 199995
$ws.Range("L135").Value = 199995
$ws.Range("N135").Value = -210135

# Row 138
$ws.Range("H138").Value = 95000
$ws.Range("J138").Value = 95000
$ws.Range("L138").Value = 95000
$ws.Range("N138").Value = -105280

# Row 139
$ws.Range("H139").Value = 99984.86
$ws.Range("J139").Value = 99984.86
$ws.Range("L139").Value = 99984.86
$ws.Range("N139").Value = -110264.86

# Row 140
$ws.Range("H140").Value = 86903
$ws.Range("I140").Value = 70709
$ws.Range("J140").Value = 95000
$ws.Range("K140").Value = 70709
$ws.Range("L140").Value = 95000
$ws.Range("M140").Value = -65529
$ws.Range("N140").Value = -105360

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 5618.7144
$ws.Range("I68").Value = 2582.75
$ws.Range("K68").Value = 2582.75
$ws.Range("M68").Value = -1833.75

# Row 71
$ws.Range("H71").Value = 5618.7144
$ws.Range("I71").Value = 2582.75
$ws.Range("K71").Value = 12913.75
$ws.Range("M71").Value = -9169.75

$ws = $wb.Worksheets.Item("WVR")
# Row 31
$ws.Range("H31").Value = 15000
$ws.Range("I31").Value = 15000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -14652
$ws.Range("N31").ClearContents()

